$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MERGED-DEFAULT+GRP2-784"
$ws.Range("A3").Value = "MERGED-DEFAULT+GRP2-152"
$ws.Range("A4").Value = "MERGED-DEFAULT+GRP2-648"
$ws.Range("A5").Value = "MERGED-GRP1-0"
$ws.Range("A6").Value = "MERGED-GRP1-808"
$ws.Range("A7").Value = "MERGED-GRP2-808"
$ws.Range("A8").Value = "MERGED-GRP2-944"
$ws.Range("A9").Value = "MERGED-GRP2-456"
$ws.Range("A10").Value = "MERGED-GRP2-904"
$ws.Range("A11").Value = "MERGED-GRP2-352"
$ws.Range("A12").Value = "MERGED-GRP2-800"
$ws.Range("A13").Value = "MERGED-GRP2-312"
